$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 1.141788666666667
$ws.Range("H2").Value2 = 3.425366
$ws.Range("I2").Value2 = 0.2235063054668694
$ws.Range("J2").Value2 = 0.2235063054668694
$ws.Range("M2").Value2 = 2.027115333333333
$ws.Range("N2").Value2 = 6.081346
$ws.Range("O2").Value2 = 0.006596284565418616
$ws.Range("P2").Value2 = 0.006596284565418615
$ws.Range("Q2").Value2 = 2.314537313626222
$ws.Range("R2").Value2 = 20.830835822636
$ws.Range("S2").Value2 = 0.001474311193024849
$ws.Range("T2").Value2 = 0.001474311193024849
$ws.Range("G3").Value2 = 1.141788666666667
$ws.Range("H3").Value2 = 3.425366
$ws.Range("I3").Value2 = 0.2235063054668694
$ws.Range("J3").Value2 = 0.2235063054668694
$ws.Range("O3").Value2 = 0.8344762556643375
$ws.Range("P3").Value2 = 0.8344762556643374
$ws.Range("Q3").Value2 = 292.8052014608062
$ws.Range("R3").Value2 = 2635.246813147256
$ws.Range("S3").Value2 = 0.1865107049033629
$ws.Range("T3").Value2 = 0.1865107049033628
$ws.Range("G4").Value2 = 1.141788666666667
$ws.Range("H4").Value2 = 3.425366
$ws.Range("I4").Value2 = 0.2235063054668694
$ws.Range("J4").Value2 = 0.2235063054668694
$ws.Range("M4").Value2 = 48.84026566666667
$ws.Range("N4").Value2 = 146.520797
$ws.Range("O4").Value2 = 0.158927459770244
$ws.Range("P4").Value2 = 0.158927459770244
$ws.Range("Q4").Value2 = 55.76526181518911
$ws.Range("R4").Value2 = 501.8873563367021
$ws.Range("S4").Value2 = 0.03552128937048177
$ws.Range("T4").Value2 = 0.03552128937048175
$ws.Range("I5").Value2 = 0.5880650598117267
$ws.Range("J5").Value2 = 0.5880650598117266
$ws.Range("M5").Value2 = 2.027115333333333
$ws.Range("N5").Value2 = 6.081346
$ws.Range("O5").Value2 = 0.006596284565418616
$ws.Range("P5").Value2 = 0.006596284565418615
$ws.Range("Q5").Value2 = 6.089754474402667
$ws.Range("R5").Value2 = 54.807790269624
$ws.Range("S5").Value2 = 0.003879044477498068
$ws.Range("T5").Value2 = 0.003879044477498067
$ws.Range("I6").Value2 = 0.5880650598117267
$ws.Range("J6").Value2 = 0.5880650598117266
$ws.Range("O6").Value2 = 0.8344762556643375
$ws.Range("P6").Value2 = 0.8344762556643374
$ws.Range("S6").Value2 = 0.4907263291987144
$ws.Range("T6").Value2 = 0.4907263291987142
$ws.Range("I7").Value2 = 0.5880650598117267
$ws.Range("J7").Value2 = 0.5880650598117266
$ws.Range("M7").Value2 = 48.84026566666667
$ws.Range("N7").Value2 = 146.520797
$ws.Range("O7").Value2 = 0.158927459770244
$ws.Range("P7").Value2 = 0.158927459770244
$ws.Range("Q7").Value2 = 146.7233864219854
$ws.Range("R7").Value2 = 1320.510477797868
$ws.Range("S7").Value2 = 0.09345968613551435
$ws.Range("T7").Value2 = 0.09345968613551429
$ws.Range("G8").Value2 = 0.9625933333333334
$ws.Range("H8").Value2 = 2.88778
$ws.Range("I8").Value2 = 0.1884286347214039
$ws.Range("J8").Value2 = 0.1884286347214039
$ws.Range("M8").Value2 = 2.027115333333333
$ws.Range("N8").Value2 = 6.081346
$ws.Range("O8").Value2 = 0.006596284565418616
$ws.Range("P8").Value2 = 0.006596284565418615
$ws.Range("Q8").Value2 = 1.951287705764444
$ws.Range("R8").Value2 = 17.56158935188
$ws.Range("S8").Value2 = 0.001242928894895699
$ws.Range("T8").Value2 = 0.001242928894895698
$ws.Range("G9").Value2 = 0.9625933333333334
$ws.Range("H9").Value2 = 2.88778
$ws.Range("I9").Value2 = 0.1884286347214039
$ws.Range("J9").Value2 = 0.1884286347214039
$ws.Range("O9").Value2 = 0.8344762556643375
$ws.Range("P9").Value2 = 0.8344762556643374
$ws.Range("Q9").Value2 = 246.8515786851644
$ws.Range("R9").Value2 = 2221.66420816648
$ws.Range("S9").Value2 = 0.1572392215622603
$ws.Range("T9").Value2 = 0.1572392215622603
$ws.Range("G10").Value2 = 0.9625933333333334
$ws.Range("H10").Value2 = 2.88778
$ws.Range("I10").Value2 = 0.1884286347214039
$ws.Range("J10").Value2 = 0.1884286347214039
$ws.Range("M10").Value2 = 48.84026566666667
$ws.Range("N10").Value2 = 146.520797
$ws.Range("O10").Value2 = 0.158927459770244
$ws.Range("P10").Value2 = 0.158927459770244
$ws.Range("Q10").Value2 = 47.01331412896224
$ws.Range("R10").Value2 = 423.1198271606601
$ws.Range("S10").Value2 = 0.02994648426424792
$ws.Range("T10").Value2 = 0.02994648426424791
